# "Generate Report for Handoff"
# A new markdown file (55684d0b-dded-44c6-ab2c-2cee52000aaf.md) has reached
# "Ready for handoff" status. Insert it as the new first data row on every
# sheet (pushing the existing 6b069158-... row down by one) and record its
# per-locale handoff artifacts.

$wb = $excel.ActiveWorkbook

$oldGuid = "6b069158-7999-440c-b13b-511123ed9da9"
$newGuid = "55684d0b-dded-44c6-ab2c-2cee52000aaf"
$oldHash = "a15f165e53887a2670d9393d054e4b8b37983a68"
$newHash = "ae07ebfffb3b4ead53edd2a81f06e2f48cbeeb96"

$mdOldUrl = "https://github.com/OpenLocalizationTest/oltest/blob/75340409c2fde23df902c8360942f64a19842b94/e2e/$oldGuid.md"
$mdNewUrl = "https://github.com/OpenLocalizationTest/oltest/blob/75340409c2fde23df902c8360942f64a19842b94/e2e/$newGuid.md"

$zhOldUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cab4477990beadad189685f2d09ed2fe7ff5af72/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$zhNewUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cab4477990beadad189685f2d09ed2fe7ff5af72/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf"

$deOldUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18445e8770b2634f3c49d56700d682da955d1960/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"
$deNewUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18445e8770b2634f3c49d56700d682da955d1960/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows(2).Insert()

# Row 2 shifted down to row 3 along with its values, but the hyperlink does
# not follow automatically - recreate it pointing at the (unchanged) old file.
$wsOverview.Range("A3").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdOldUrl, "", "", "$oldGuid.md")

# New row 2: the freshly handed-off file.
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-29-20 04:29:20"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdNewUrl, "", "", "$newGuid.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows(2).Insert()

$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdOldUrl, "", "", "$oldGuid.md")
$wsZhCn.Range("B3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $mdOldUrl, "", "", ".md")
$wsZhCn.Range("D3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhOldUrl, "", "", "$oldGuid.$oldHash.zh-cn.xlf")
$wsZhCn.Range("E3").NumberFormat = $dateFmt

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-20 04:29:17"
$wsZhCn.Range("E2").NumberFormat = $dateFmt
$wsZhCn.Range("H2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I2").Value = "Include"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdNewUrl, "", "", "$newGuid.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $mdNewUrl, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhNewUrl, "", "", "$newGuid.$newHash.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn" with de-de handoff artifacts.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows(2).Insert()

$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdOldUrl, "", "", "$oldGuid.md")
$wsDeDe.Range("B3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $mdOldUrl, "", "", ".md")
$wsDeDe.Range("D3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deOldUrl, "", "", "$oldGuid.$oldHash.de-de.xlf")
$wsDeDe.Range("E3").NumberFormat = $dateFmt

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-20 04:29:20"
$wsDeDe.Range("E2").NumberFormat = $dateFmt
$wsDeDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I2").Value = "Include"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdNewUrl, "", "", "$newGuid.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $mdNewUrl, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deNewUrl, "", "", "$newGuid.$newHash.de-de.xlf")
